# Project DesignFirst: update the "From" value for rule R20 (row 10) in the
# Rules sheet from 18 to 100.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 100
